# Apply the "Updated cryptos list" data refresh described by the commit diff.
# For cells whose new text would otherwise be auto-parsed by Excel as a number
# (e.g. "1.007", "19.90"), a leading apostrophe forces text entry; the style is
# then reset to "Normal" so no stray quote-prefix style survives in the saved file.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.132.83'
$ws.Range("E2").Value = '  -4.35%  '
$ws.Range("D3").Value = '1.654.33'
$ws.Range("E3").Value = '  -3.30%  '
$ws.Range("E4").Value = '  +0.29%  '
$ws.Range("D5").Value = '''215.42'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.95%  '
$ws.Range("D6").Value = '''0.5096'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.77%  '
$ws.Range("E7").Value = '  +0.30%  '
$ws.Range("D8").Value = '''0.2578'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.15%  '
$ws.Range("D9").Value = '''0.06404'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.39%  '
$ws.Range("D10").Value = '''19.90'
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = '''0.07798'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.69%  '
$ws.Range("D12").Value = '1.658.52'
$ws.Range("E12").Value = '  -2.97%  '
$ws.Range("D13").Value = '''4.279'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -5.07%  '
$ws.Range("D14").Value = '1.882.28'
$ws.Range("E14").Value = '  -3.30%  '
$ws.Range("D15").Value = '''0.5508'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -5.30%  '
$ws.Range("D16").Value = '0.0₅8004'
$ws.Range("E16").Value = '  -2.81%  '
$ws.Range("D17").Value = '''63.93'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -6.27%  '
$ws.Range("D18").Value = '26.155.24'
$ws.Range("E18").Value = '  -4.31%  '
$ws.Range("E19").Value = '  +0.26%  '
$ws.Range("D20").Value = '''209.17'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -7.35%  '
$ws.Range("D21").Value = '''4.404'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.83%  '
$ws.Range("E22").Value = '  -3.28%  '
$ws.Range("D23").Value = '''6.028'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.37%  '
$ws.Range("D24").Value = '''1.007'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.30%  '
$ws.Range("E25").Value = '  -0.58%  '
$ws.Range("D26").Value = '''1.740'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.58%  '
$ws.Range("D27").Value = '''0.1177'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.40%  '
$ws.Range("D28").Value = '''6.973'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.57%  '
$ws.Range("D29").Value = '''15.81'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.98%  '
$ws.Range("D30").Value = '''0.05103'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.82%  '
$ws.Range("D31").Value = '''1.242'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.77%  '
$ws.Range("E32").Value = '  -4.04%  '
$ws.Range("D33").Value = '''3.219'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -6.25%  '
$ws.Range("D34").Value = '''1.567'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.33%  '
$ws.Range("D35").Value = '''2.754'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.03%  '
$ws.Range("D36").Value = '''0.9277'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.33%  '
$ws.Range("E37").Value = '  -1.23%  '
$ws.Range("B38").Value = 'Maker'
$ws.Range("C38").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D38").Value = '1.164.20'
$ws.Range("E38").Value = '  +7.74%  '
$ws.Range("B39").Value = 'ImmutableX'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D39").Value = '''0.5687'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.37%  '
$ws.Range("E40").Value = '  -2.84%  '
$ws.Range("E41").Value = '  +0.26%  '
$ws.Range("D42").Value = '''2.555'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.27%  '
$ws.Range("D43").Value = '''0.8333'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.72%  '
$ws.Range("D44").Value = '''5.648'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.35%  '
$ws.Range("D45").Value = '''100.39'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.36%  '
$ws.Range("D46").Value = '1.792.66'
$ws.Range("E46").Value = '  -3.30%  '
$ws.Range("E47").Value = '  +0.09%  '
$ws.Range("D48").Value = '''0.4550'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.49%  '
$ws.Range("D49").Value = '''55.71'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.48%  '
$ws.Range("D50").Value = '''1.008'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.57%  '
$ws.Range("D51").Value = '''7.843'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.84%  '
